$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Force column D (Price) to remain plain text so values like "12.40" or
# "0.0000317" are not auto-coerced into Double/scientific-notation by Excel's
# native value parser.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "70.820.12"
$ws.Cells.Item(2, 5).Value = "  +0.89%  "

$ws.Cells.Item(3, 4).Value = "3.586.34"
$ws.Cells.Item(3, 5).Value = "  +0.08%  "

$ws.Cells.Item(4, 5).Value = "  +0.12%  "

$ws.Cells.Item(5, 4).Value = "586.49"
$ws.Cells.Item(5, 5).Value = "  +1.33%  "

$ws.Cells.Item(6, 4).Value = "186.44"
$ws.Cells.Item(6, 5).Value = "  -0.22%  "

$ws.Cells.Item(7, 4).Value = "3.573.60"
$ws.Cells.Item(7, 5).Value = "  -0.15%  "

$ws.Cells.Item(8, 5).Value = "  +0.45%  "

$ws.Cells.Item(9, 5).Value = "  +0.03%  "

$ws.Cells.Item(10, 4).Value = "0.212"
$ws.Cells.Item(10, 5).Value = "  +15.55%  "

$ws.Cells.Item(11, 5).Value = "  -0.08%  "

$ws.Cells.Item(12, 4).Value = "54.32"
$ws.Cells.Item(12, 5).Value = "  -1.61%  "

$ws.Cells.Item(13, 4).Value = "0.0000317"
$ws.Cells.Item(13, 5).Value = "  +3.56%  "

$ws.Cells.Item(14, 5).Value = "  -0.13%  "

$ws.Cells.Item(15, 4).Value = "4.155.69"
$ws.Cells.Item(15, 5).Value = "  -0.07%  "

$ws.Cells.Item(16, 4).Value = "19.58"
$ws.Cells.Item(16, 5).Value = "  -0.68%  "

$ws.Cells.Item(17, 4).Value = "70.828.65"
$ws.Cells.Item(17, 5).Value = "  +1.07%  "

$ws.Cells.Item(18, 4).Value = "3.575.35"
$ws.Cells.Item(18, 5).Value = "  +0.00%  "

$ws.Cells.Item(19, 2).Value = "BitcoinCash"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(19, 4).Value = "575.31"
$ws.Cells.Item(19, 5).Value = "  +16.52%  "

$ws.Cells.Item(20, 2).Value = "Uniswap"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Cells.Item(20, 4).Value = "12.41"
$ws.Cells.Item(20, 5).Value = "  -1.49%  "

$ws.Cells.Item(21, 5).Value = "  -0.17%  "

$ws.Cells.Item(22, 5).Value = "  -2.19%  "

$ws.Cells.Item(23, 4).Value = "17.79"
$ws.Cells.Item(23, 5).Value = "  -7.92%  "

$ws.Cells.Item(24, 4).Value = "4.64"
$ws.Cells.Item(24, 5).Value = "  +5.55%  "

$ws.Cells.Item(25, 4).Value = "4.91"
$ws.Cells.Item(25, 5).Value = "  -1.16%  "

$ws.Cells.Item(26, 4).Value = "95.43"
$ws.Cells.Item(26, 5).Value = "  -1.34%  "

$ws.Cells.Item(27, 4).Value = "11.47"
$ws.Cells.Item(27, 5).Value = "  -0.46%  "

$ws.Cells.Item(28, 4).Value = "2.95"
$ws.Cells.Item(28, 5).Value = "  -0.35%  "

$ws.Cells.Item(29, 4).Value = "9.16"
$ws.Cells.Item(29, 5).Value = "  -2.11%  "

$ws.Cells.Item(30, 4).Value = "32.21"
$ws.Cells.Item(30, 5).Value = "  +1.51%  "

$ws.Cells.Item(31, 4).Value = "7.32"
$ws.Cells.Item(31, 5).Value = "  -5.68%  "

$ws.Cells.Item(32, 4).Value = "12.40"
$ws.Cells.Item(32, 5).Value = "  +2.12%  "

$ws.Cells.Item(33, 4).Value = "65.02"
$ws.Cells.Item(33, 5).Value = "  -1.32%  "

$ws.Cells.Item(34, 4).Value = "0.115"
$ws.Cells.Item(34, 5).Value = "  -0.85%  "

$ws.Cells.Item(35, 2).Value = "Fetch.AI"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(35, 4).Value = "3.33"
$ws.Cells.Item(35, 5).Value = "  +2.88%  "

$ws.Cells.Item(36, 2).Value = "Bittensor"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(36, 4).Value = "563.93"
$ws.Cells.Item(36, 5).Value = "  -1.96%  "

$ws.Cells.Item(37, 4).Value = "0.416"
$ws.Cells.Item(37, 5).Value = "  +0.56%  "

$ws.Cells.Item(38, 4).Value = "0.0₃0804"
$ws.Cells.Item(38, 5).Value = "  +1.02%  "

$ws.Cells.Item(39, 4).Value = "37.80"
$ws.Cells.Item(39, 5).Value = "  -2.83%  "

$ws.Cells.Item(40, 5).Value = "  +0.02%  "

$ws.Cells.Item(41, 4).Value = "3.394.15"
$ws.Cells.Item(41, 5).Value = "  +6.34%  "

$ws.Cells.Item(42, 4).Value = "3.12"
$ws.Cells.Item(42, 5).Value = "  -1.45%  "

$ws.Cells.Item(43, 5).Value = "  -0.52%  "

$ws.Cells.Item(44, 4).Value = "3.38"
$ws.Cells.Item(44, 5).Value = "  -3.16%  "

$ws.Cells.Item(45, 4).Value = "3.58"
$ws.Cells.Item(45, 5).Value = "  -1.74%  "

$ws.Cells.Item(46, 4).Value = "0.0447"
$ws.Cells.Item(46, 5).Value = "  +0.94%  "

$ws.Cells.Item(47, 4).Value = "2.97"
$ws.Cells.Item(47, 5).Value = "  -3.27%  "

$ws.Cells.Item(48, 4).Value = "9.40"
$ws.Cells.Item(48, 5).Value = "  -1.01%  "

$ws.Cells.Item(49, 5).Value = "  +0.64%  "

$ws.Cells.Item(50, 4).Value = "1.00"
$ws.Cells.Item(50, 5).Value = "  +0.14%  "

$ws.Cells.Item(51, 4).Value = "1.42"
$ws.Cells.Item(51, 5).Value = "  -9.07%  "
